$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C header — copy the header style (bold, centered, bordered) from B1, then set text
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Coord: normal vector scan"

# Updated B values (column B, rows 2-15)
$bValues = @(
    0.2992193909499149,
    1.154542479549355,
    0.2417729669738042,
    0.5642115139385254,
    1.634349092659157,
    0.2590600675317279,
    0.5813071460277248,
    0.9783082765882792,
    1.8147533940247,
    1.63659575012755,
    0.5823097362680166,
    1.530902043419413,
    2.888176878425738,
    1.359093703802407
)

# New column C values (inline strings of normal vectors)
$cValues = @(
    "[0.         0.32002289 0.94740981]",
    "[-0.44071658  0.52582622  0.72751336]",
    "[-0.00421424  0.01625641  0.99985897]",
    "[-6.33044593e-04  2.85055526e-01 -9.58510796e-01]",
    "[0.72484588 0.28839264 0.62564218]",
    "[-0.7390617  -0.27315406  0.6157716 ]",
    "[0.         0.31535459 0.94897391]",
    "[ 0.         -0.30877167  0.95113619]",
    "[-0.72754179  0.25842599  0.63553045]",
    "[ 0.72923746 -0.27975851  0.62445809]",
    "[ 6.40813527e-04 -2.84752602e-01 -9.58600827e-01]",
    "[-0.72573589 -0.28922681  0.62422373]",
    "[0.71129802 0.26422727 0.65133637]",
    "[ 0.         -0.30244368  0.95316726]"
)

for ($i = 0; $i -lt 14; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}
